$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.827.87"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "3.643.61"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  +17.89%  "
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'225.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.78%  "
$ws.Range("D7").Value = "'647.92"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("E8").Value = "  -3.54%  "
$ws.Range("E9").Value = "  +4.44%  "
$ws.Range("D10").Value = "'0.999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "3.638.99"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "'52.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.02%  "
$ws.Range("D13").Value = "'0.220"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.38%  "
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "'6.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.86%  "
$ws.Range("D16").Value = "4.323.11"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "'24.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +32.25%  "
$ws.Range("D18").Value = "95.483.15"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "'9.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "'13.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.07%  "
$ws.Range("D21").Value = "3.637.04"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "'0.310"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +48.43%  "
$ws.Range("D23").Value = "'0.538"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.57%  "
$ws.Range("D24").Value = "'538.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").Value = "'133.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.91%  "
$ws.Range("D26").Value = "'3.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").Value = "'7.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.76%  "
$ws.Range("E28").Value = "  -8.25%  "
$ws.Range("D29").Value = "'13.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").Value = "3.810.67"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "'13.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.51%  "
$ws.Range("E32").Value = "  +5.95%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'1.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.78%  "
$ws.Range("D35").Value = "'0.649"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.65%  "
$ws.Range("D36").Value = "'33.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "'0.0571"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +23.86%  "
$ws.Range("D40").Value = "'8.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "'606.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'7.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.17%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.92%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.503"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").Value = "'41.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E48").Value = "  -6.37%  "
$ws.Range("D49").Value = "'9.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.63%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'235.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.62%  "
